# Add a "Total Ukr Cases" column (I) with per-city case totals, and update
# the sheet's header, selection and column width accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lookup table: City -> Total Ukr Cases value
$caseLookup = @{
    "Odessa"                 = 323448
    "Bilhorod-Dnistrovskyi"  = 241629
    "Kherson"                = 311484
    "Lviv"                   = 371768
    "Ivano-Frankivsk"        = 344581
    "Drohobych"              = 369162
    "Kyiv"                   = 411341
    "Bila Tserkva"           = 416328
    "Lubny"                  = 406495
    "Kharkiv"                = 407589
    "Kramatorsk"             = 404793
    "Izyum"                  = 392189
}

# New column header
$ws.Cells.Item(1, 9).Value = "Total Ukr Cases"

# Fill in the new column using each row's City (column A) to look up the value
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $city = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 9).Value = $caseLookup[$city]
}

# Widen column H to fit its content, matching the post-edit layout
$ws.Columns.Item(8).ColumnWidth = 11.27

# Update the active selection to reflect where the author left off editing
$ws.Range("G10").Select() | Out-Null
